$d = $word.ActiveDocument

# --- Change 1: split the title run into a plain part + a new bold part ---
$p1 = $d.Paragraphs.Item(1)
$titleStart = $p1.Range.Start
$titleEnd = $p1.Range.End
$splitAt = $titleStart + 22   # length of "SEVES - Alimentaire - "

$rTail = $d.Range($splitAt, $titleEnd)
$rTail.Text = "Investigation cas humain"
$rTail.Font.Bold = 1

# --- Change 2: move "{%p endfor %}" text to a new paragraph before the
#     section-break paragraph, leaving that paragraph's run text-less ---
# (there are two "{%p endfor %}" paragraphs in the doc; the one we want is
#  the one that carries a section break - i.e. its paragraph mark is the
#  "end of section" mark rather than a plain paragraph mark.
#  NB: Paragraph.Index is unreliable in this host, so the 1-based position
#  is captured from the enumeration instead of read back afterwards.)
$targetPos = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t.Length -eq 14 -and $t.Substring(0, 13) -eq "{%p endfor %}") {
        $markCode = [int][char]$t[13]
        if ($markCode -eq 12) {
            $targetPos = $i
        }
    }
}

$target = $d.Paragraphs.Item($targetPos)
$target.Range.InsertParagraphBefore()

$pNew = $d.Paragraphs.Item($targetPos)
$pOld = $d.Paragraphs.Item($targetPos + 1)

$pNew.Range.Text = "{%p endfor %}"

$oldStart = $pOld.Range.Start
$oldTextRange = $d.Range($oldStart, $oldStart + 13)
$oldTextRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr/></w:r></w:p>")
